$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"1.046867666666667"
$ws.Range("H2").Value = [double]"3.140603"
$ws.Range("I2").Value = [double]"0.000687505225377314"
$ws.Range("J2").Value = [double]"0.000687505225377314"
$ws.Range("M2").Value = [double]"547.450775"
$ws.Range("N2").Value = [double]"1642.352325"
$ws.Range("O2").Value = [double]"0.8253533007282613"
$ws.Range("P2").Value = [double]"0.8253533007282614"
$ws.Range("Q2").Value = [double]"573.1085154391084"
$ws.Range("R2").Value = [double]"5157.976638951975"
$ws.Range("S2").Value = [double]"0.0005674347070330932"
$ws.Range("T2").Value = [double]"0.0005674347070330933"

# Row 3
$ws.Range("G3").Value = [double]"1.046867666666667"
$ws.Range("H3").Value = [double]"3.140603"
$ws.Range("I3").Value = [double]"0.000687505225377314"
$ws.Range("J3").Value = [double]"0.000687505225377314"
$ws.Range("O3").Value = [double]"0.002183077622430991"
$ws.Range("P3").Value = [double]"0.002183077622430991"
$ws.Range("Q3").Value = [double]"1.515884620774889"
$ws.Range("R3").Value = [double]"13.642961586974"
$ws.Range("S3").Value = [double]"1.500877272825589E-06"
$ws.Range("T3").Value = [double]"1.500877272825589E-06"

# Row 4
$ws.Range("G4").Value = [double]"1.046867666666667"
$ws.Range("H4").Value = [double]"3.140603"
$ws.Range("I4").Value = [double]"0.000687505225377314"
$ws.Range("J4").Value = [double]"0.000687505225377314"
$ws.Range("M4").Value = [double]"114.393852"
$ws.Range("N4").Value = [double]"343.181556"
$ws.Range("O4").Value = [double]"0.1724636216493076"
$ws.Range("P4").Value = [double]"0.1724636216493076"
$ws.Range("Q4").Value = [double]"119.755224924252"
$ws.Range("R4").Value = [double]"1077.797024318268"
$ws.Range("S4").Value = [double]"0.000118569641071395"
$ws.Range("T4").Value = [double]"0.000118569641071395"

# Row 5
$ws.Range("H5").Value = [double]"4442.55542"
$ws.Range("I5").Value = [double]"0.9725138978974124"
$ws.Range("J5").Value = [double]"0.9725138978974125"
$ws.Range("M5").Value = [double]"547.450775"
$ws.Range("N5").Value = [double]"1642.352325"
$ws.Range("O5").Value = [double]"0.8253533007282613"
$ws.Range("P5").Value = [double]"0.8253533007282614"
$ws.Range("Q5").Value = [double]"810693.4692198168"
$ws.Range("R5").Value = [double]"7296241.222978352"
$ws.Range("S5").Value = [double]"0.8026675556337366"
$ws.Range("T5").Value = [double]"0.8026675556337368"

# Row 6
$ws.Range("H6").Value = [double]"4442.55542"
$ws.Range("I6").Value = [double]"0.9725138978974124"
$ws.Range("J6").Value = [double]"0.9725138978974125"
$ws.Range("O6").Value = [double]"0.002183077622430991"
$ws.Range("P6").Value = [double]"0.002183077622430991"
$ws.Range("S6").Value = [double]"0.002123073328002978"
$ws.Range("T6").Value = [double]"0.002123073328002979"

# Row 7
$ws.Range("H7").Value = [double]"4442.55542"
$ws.Range("I7").Value = [double]"0.9725138978974124"
$ws.Range("J7").Value = [double]"0.9725138978974125"
$ws.Range("M7").Value = [double]"114.393852"
$ws.Range("N7").Value = [double]"343.181556"
$ws.Range("O7").Value = [double]"0.1724636216493076"
$ws.Range("P7").Value = [double]"0.1724636216493076"
$ws.Range("Q7").Value = [double]"169400.3424057593"
$ws.Range("R7").Value = [double]"1524603.081651833"
$ws.Range("S7").Value = [double]"0.1677232689356727"
$ws.Range("T7").Value = [double]"0.1677232689356727"

# Row 8
$ws.Range("G8").Value = [double]"40.80635833333333"
$ws.Range("H8").Value = [double]"122.419075"
$ws.Range("I8").Value = [double]"0.02679859687721029"
$ws.Range("J8").Value = [double]"0.0267985968772103"
$ws.Range("M8").Value = [double]"547.450775"
$ws.Range("N8").Value = [double]"1642.352325"
$ws.Range("O8").Value = [double]"0.8253533007282613"
$ws.Range("P8").Value = [double]"0.8253533007282614"
$ws.Range("Q8").Value = [double]"22339.47249451104"
$ws.Range("R8").Value = [double]"201055.2524505994"
$ws.Range("S8").Value = [double]"0.02211831038749159"
$ws.Range("T8").Value = [double]"0.0221183103874916"

# Row 9
$ws.Range("G9").Value = [double]"40.80635833333333"
$ws.Range("H9").Value = [double]"122.419075"
$ws.Range("I9").Value = [double]"0.02679859687721029"
$ws.Range("J9").Value = [double]"0.0267985968772103"
$ws.Range("O9").Value = [double]"0.002183077622430991"
$ws.Range("P9").Value = [double]"0.002183077622430991"
$ws.Range("Q9").Value = [double]"59.08839578959444"
$ws.Range("R9").Value = [double]"531.79556210635"
$ws.Range("S9").Value = [double]"5.850341715518682E-05"
$ws.Range("T9").Value = [double]"5.850341715518684E-05"

# Row 10
$ws.Range("G10").Value = [double]"40.80635833333333"
$ws.Range("H10").Value = [double]"122.419075"
$ws.Range("I10").Value = [double]"0.02679859687721029"
$ws.Range("J10").Value = [double]"0.0267985968772103"
$ws.Range("M10").Value = [double]"114.393852"
$ws.Range("N10").Value = [double]"343.181556"
$ws.Range("O10").Value = [double]"0.1724636216493076"
$ws.Range("P10").Value = [double]"0.1724636216493076"
$ws.Range("Q10").Value = [double]"4667.996515842299"
$ws.Range("R10").Value = [double]"42011.96864258069"
$ws.Range("S10").Value = [double]"0.004621783072563512"
$ws.Range("T10").Value = [double]"0.004621783072563513"
